$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Status" values in column C for rows 19-22 and 25
$ws.Range("C19").Value = "DONE"
$ws.Range("C20").Value = "TODO"
$ws.Range("C21").Value = "Not now"
$ws.Range("C22").Value = "Not now"
$ws.Range("C25").Value = "Works Unique, does not duplicate"

# Row 20 wraps across one more line now that its text records the dropped
# item, so its autofit height grows from 78.75 to 105 points.
$ws.Rows("20").RowHeight = 105

# Update the selected/visible range to reflect scrolling down to the new rows
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C27").Select()
